$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'" + '26.957.41'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.31%  '
$ws.Range('D3').Value = "'" + '1.550.67'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.82%  '
$ws.Range('E4').Value = '  +0.27%  '
$ws.Range('D5').Value = "'" + '206.10'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.16%  '
$ws.Range('D6').Value = "'" + '0.485'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.33%  '
$ws.Range('E7').Value = '  +0.17%  '
$ws.Range('E8').Value = '  +0.45%  '
$ws.Range('D9').Value = "'" + '21.48'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.46%  '
$ws.Range('D10').Value = "'" + '0.0583'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.61%  '
$ws.Range('D11').Value = "'" + '0.0858'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.82%  '
$ws.Range('D12').Value = "'" + '1.771.34'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.87%  '
$ws.Range('D13').Value = "'" + '1.551.05'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.71%  '
$ws.Range('E14').Value = '  -0.69%  '
$ws.Range('E15').Value = '  -0.48%  '
$ws.Range('D16').Value = "'" + '26.949.65'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.28%  '
$ws.Range('D17').Value = "'" + '61.60'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.52%  '
$ws.Range('D18').Value = "'" + '215.02'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.14%  '
$ws.Range('D19').Value = "'" + '0.0₃0685'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.41%  '
$ws.Range('D20').Value = "'" + '7.23'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -2.06%  '
$ws.Range('E21').Value = '  +0.23%  '
$ws.Range('E22').Value = '  -2.46%  '
$ws.Range('D23').Value = "'" + '9.20'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.01%  '
$ws.Range('D24').Value = "'" + '1.95'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -3.12%  '
$ws.Range('D25').Value = "'" + '153.41'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.51%  '
$ws.Range('E26').Value = '  -0.46%  '
$ws.Range('D27').Value = "'" + '14.84'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.67%  '
$ws.Range('E28').Value = '  +0.13%  '
$ws.Range('E29').Value = '  -0.03%  '
$ws.Range('E30').Value = '  -1.11%  '
$ws.Range('E31').Value = '  -1.08%  '
$ws.Range('D33').Value = "'" + '1.371.20'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.53%  '
$ws.Range('D34').Value = "'" + '2.96'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.20%  '
$ws.Range('E35').Value = '  +1.56%  '
$ws.Range('E36').Value = '  +4.62%  '
$ws.Range('E37').Value = '  +0.18%  '
$ws.Range('E38').Value = '  +0.22%  '
$ws.Range('D39').Value = "'" + '0.519'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -2.01%  '
$ws.Range('D40').Value = "'" + '0.807'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.67%  '
$ws.Range('E41').Value = '  +0.21%  '
$ws.Range('B42').Value = 'MXToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D42').Value = "'" + '2.30'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +5.13%  '
$ws.Range('B43').Value = 'WEMIXToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D43').Value = "'" + '0.985'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.70%  '
$ws.Range('B44').Value = 'FraxShare'
$ws.Range('C44').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D44').Value = "'" + '5.48'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.89%  '
$ws.Range('D45').Value = "'" + '63.70'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.01%  '
$ws.Range('D46').Value = "'" + '1.76'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.62%  '
$ws.Range('B47').Value = 'RocketPoolETH'
$ws.Range('C47').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D47').Value = "'" + '1.685.17'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.00%  '
$ws.Range('B48').Value = 'mCoin'
$ws.Range('C48').Value = 'https://coinranking.com/coin/fzVgyjBcRc9+mcoin-mcoin'
$ws.Range('D48').Value = "'" + '2.25'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -3.03%  '
$ws.Range('D49').Value = "'" + '86.30'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.54%  '
$ws.Range('D50').Value = "'" + '0.0508'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.99%  '
$ws.Range('D51').Value = "'" + '0.0953'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.01%  '
